# Updates the cryptos list (Price / Volume(1h) columns) for the latest
# GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (column D) / Volume(1h) (column E) values.
# DNum flags a Price value that looks like a plain decimal number (e.g.
# "1.00", "586.51"): Excel would otherwise auto-convert it to a numeric
# cell, so those get the column pre-formatted as Text before the value is
# written, keeping them as text just like the source feed exports them.
$updates = @(
    @{Row=2; D='63.122.17'; E='  -0.05%  '},
    @{Row=3; D='3.054.11'; E='  -0.10%  '},
    @{Row=4; D='1.00'; DNum=$true; E='  +0.01%  '},
    @{Row=5; D='586.51'; DNum=$true; E='  -0.51%  '},
    @{Row=6; D='151.79'; DNum=$true; E='  -0.93%  '},
    @{Row=7; E='  +0.04%  '},
    @{Row=8; D='0.537'; DNum=$true; E='  -1.77%  '},
    @{Row=9; D='3.053.93'; E='  -0.38%  '},
    @{Row=10; E='  -2.65%  '},
    @{Row=11; D='5.87'; DNum=$true; E='  +0.57%  '},
    @{Row=12; E='  -2.60%  '},
    @{Row=13; E='  -2.25%  '},
    @{Row=14; D='36.27'; DNum=$true; E='  -2.78%  '},
    @{Row=15; E='  +1.94%  '},
    @{Row=16; D='3.557.18'; E='  -0.18%  '},
    @{Row=17; E='  -0.96%  '},
    @{Row=18; D='63.133.62'; E='  -0.17%  '},
    @{Row=19; D='3.052.99'; E='  -0.34%  '},
    @{Row=20; D='477.66'; DNum=$true; E='  +0.35%  '},
    @{Row=21; E='  -2.46%  '},
    @{Row=22; E='  -1.42%  '},
    @{Row=23; D='7.52'; DNum=$true; E='  -0.18%  '},
    @{Row=24; D='2.41'; DNum=$true; E='  +0.75%  '},
    @{Row=25; D='82.25'; DNum=$true; E='  +1.50%  '},
    @{Row=26; D='12.72'; DNum=$true; E='  -1.68%  '},
    @{Row=27; D='10.52'; DNum=$true; E='  +5.10%  '},
    @{Row=28; E='  +0.36%  '},
    @{Row=29; D='7.36'; DNum=$true; E='  +0.66%  '},
    @{Row=31; E='  +0.01%  '},
    @{Row=32; D='2.20'; DNum=$true; E='  +0.20%  '},
    @{Row=33; D='27.62'; DNum=$true; E='  +1.59%  '},
    @{Row=34; E='  -2.39%  '},
    @{Row=35; E='  +1.52%  '},
    @{Row=36; D='0.0₃0820'; E='  -3.04%  '},
    @{Row=37; D='3.26'; DNum=$true; E='  -2.51%  '},
    @{Row=38; E='  -2.99%  '},
    @{Row=39; E='  -0.19%  '},
    @{Row=40; E='  -0.81%  '},
    @{Row=41; D='50.44'; DNum=$true; E='  -0.16%  '},
    @{Row=42; D='434.48'; DNum=$true; E='  -2.02%  '},
    @{Row=43; D='0.289'; DNum=$true; E='  +1.44%  '},
    @{Row=44; D='0.115'; DNum=$true; E='  +2.81%  '},
    @{Row=45; D='0.0361'; DNum=$true; E='  -0.20%  '},
    @{Row=46; D='2.829.10'; E='  +1.17%  '},
    @{Row=47; D='38.26'; DNum=$true; E='  -4.39%  '},
    @{Row=48; D='128.53'; DNum=$true; E='  -2.59%  '},
    @{Row=49; E='  -0.01%  '},
    @{Row=50; D='25.10'; DNum=$true; E='  -0.29%  '},
    @{Row=51; D='0.110'; DNum=$true; E='  -1.35%  '}
)

foreach ($u in $updates) {
    if ($u.ContainsKey('D')) {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($u.ContainsKey('DNum')) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($u.ContainsKey('E')) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
